# Generate Report for Handoff
# Moves the localization status from "In Translation" to "Ready for handoff"
# and refreshes the handoff timestamps on the Overview sheet as well as the
# per-language (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-28-19 10:28:37"

# --- zh-cn detail sheet ----------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-19 10:28:34"

# --- de-de detail sheet ----------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-19 10:28:37"
